$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "(0.47)"
$ws.Range("C4").Value = "(0.67)"
$ws.Range("D4").Value = "(0.19)"
$ws.Range("E4").Value = "(0.66)"
$ws.Range("F4").Value = "(0.05)"
$ws.Range("G4").Value = "(0.93)"
$ws.Range("H4").Value = "(0.35)"
$ws.Range("I4").Value = "(0.96)"
$ws.Range("J4").Value = "(0.57)"

$ws.Range("B6").Value = "(0.33)"
$ws.Range("C6").Value = "(0.3)"
$ws.Range("D6").Value = "(0.31)"
$ws.Range("E6").Value = "(1.24)"
$ws.Range("F6").Value = "(0.23)"
$ws.Range("G6").Value = "(0.22)"
$ws.Range("H6").Value = "(0.13)"
$ws.Range("I6").Value = "(0.16)"
$ws.Range("J6").Value = "(1.27)"
